$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the npv value for the "Current Policies" scenario (row 3, column B):
# it was mistakenly "$1552 Billion" and should be "$1467 Billion".
$ws.Range("B3").Value = "$1467 Billion"

# Update the active cell / selection to B4 (matches corrected sheet view state).
$ws.Activate()
$ws.Range("B4").Select()
